$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: 2026-02-11 (serial 46064), 四方坪站
$ws.Range("A22").Value = 46064
$ws.Range("B22").Value = "四方坪站"
$ws.Range("C22").Value = 12391.96
$ws.Range("D22").Value = 11433.34
$ws.Range("E22").Value = 4602.46
$ws.Range("F22").Value = 508

# Row 23: 2026-02-11 (serial 46064), 高岭站
$ws.Range("A23").Value = 46064
$ws.Range("B23").Value = "高岭站"
$ws.Range("C23").Value = 4219.9399999999996
$ws.Range("D23").Value = 3845.37
$ws.Range("E23").Value = 1206.98
$ws.Range("F23").Value = 147

# Row 24: 2026-02-12 (serial 46065), 四方坪站
$ws.Range("A24").Value = 46065
$ws.Range("B24").Value = "四方坪站"
$ws.Range("C24").Value = 10799.29
$ws.Range("D24").Value = 10026.299999999999
$ws.Range("E24").Value = 4016.22
$ws.Range("F24").Value = 461

# Row 25: 2026-02-12 (serial 46065), 高岭站
$ws.Range("A25").Value = 46065
$ws.Range("B25").Value = "高岭站"
$ws.Range("C25").Value = 3773.31
$ws.Range("D25").Value = 3361.19
$ws.Range("E25").Value = 1098.71
$ws.Range("F25").Value = 119

# Match the author's final selection state
$ws.Range("J23").Select() | Out-Null
